$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column by column to reproduce the original shared-string insertion order:
# Stone, Gem, Material, Guh, Shiny
$ws.Range("A7").Value = "Stone"
$ws.Range("A8").Value = "Gem"

$ws.Range("C7").Value = "None"
$ws.Range("C8").Value = "None"

$ws.Range("E7").Value = "Material"
$ws.Range("E8").Value = "Material"

$ws.Range("B7").Value = "Guh"
$ws.Range("B8").Value = "Shiny"

$ws.Range("D7").Value = 1
$ws.Range("D8").Value = 1

# Apply the same style (wrap text) as the rest of the table
$ws.Range("A7:E8").WrapText = $true

# Update the selection to match the post-edit state
$ws.Range("B9").Select()
